$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Every data row (2..416) has its "Förändrad" (column C) date bumped from
#    45192 to 45202. Do this first for the whole sheet.
# ---------------------------------------------------------------------------
for ($r = 2; $r -le 416; $r++) {
    $ws.Cells.Item($r, 3).Value = 45202
}

# ---------------------------------------------------------------------------
# 2) Rows 25-27 got re-ordered (new entries are sorted to the top of that
#    block). The row that used to be #27 (A 23018-2023) becomes #25 (and
#    picks up a couple of updated stats + one extra species), the old #25
#    (A 51520-2020) becomes #26, and the old #26 (A 65342-2020) becomes #27.
# ---------------------------------------------------------------------------

# Snapshot the plain (non-formula) values of the old rows 25 and 26 before
# anything is overwritten, then shift them down one slot.
$row25vals = $ws.Range("A25:R25").Value()
$row26vals = $ws.Range("A26:R26").Value()

$ws.Range("A26:R26").Value = $row25vals
$ws.Range("A27:R27").Value = $row26vals

# Row 25 becomes the old "A 23018-2023" entry, with updated figures.
$ws.Cells.Item(25, 1).Value = "A 23018-2023"
$ws.Cells.Item(25, 2).Value = 45073
$ws.Cells.Item(25, 3).Value = 45202
$ws.Cells.Item(25, 4).Value = "NORRBOTTENS LÄN"
$ws.Cells.Item(25, 5).Value = "JOKKMOKK"
$ws.Cells.Item(25, 6).Value = "SCA"
$ws.Cells.Item(25, 7).Value = 2.9
$ws.Cells.Item(25, 8).Value = 3
$ws.Cells.Item(25, 9).Value = 1
$ws.Cells.Item(25, 10).Value = 6
$ws.Cells.Item(25, 11).Value = 0
$ws.Cells.Item(25, 12).Value = 0
$ws.Cells.Item(25, 13).Value = 0
$ws.Cells.Item(25, 14).Value = 0
$ws.Cells.Item(25, 15).Value = 6
$ws.Cells.Item(25, 16).Value = 0
$ws.Cells.Item(25, 17).Value = 8
$ws.Range("R25").Value = "Järpe`r`nKnottrig blåslav`r`nMotaggsvamp`r`nOrange taggsvamp`r`nTalltita`r`nVedskivlav`r`nNorrlandslav`r`nRevlummer"

# Fix up the C (Förändrad) column on the shifted rows too (the loop above
# already set them, but do it again defensively in case values got clobbered
# by the block-copy above).
$ws.Cells.Item(26, 3).Value = 45202
$ws.Cells.Item(27, 3).Value = 45202

# ---------------------------------------------------------------------------
# 3) The HYPERLINK formulas in S/T/V/W/X/Y embed the "Beteckning" text
#    directly, so they must be rewritten per destination row to match
#    whichever entry now lives there.
# ---------------------------------------------------------------------------
function Set-BeteckningLinks($rowNum, $beteckning) {
    $ws.Range("S$rowNum").Formula = '=HYPERLINK("https://klasma.github.io/Logging_JOKKMOKK/artfynd/' + $beteckning + '.xlsx", "' + $beteckning + '")'
    $ws.Range("T$rowNum").Formula = '=HYPERLINK("https://klasma.github.io/Logging_JOKKMOKK/kartor/' + $beteckning + '.png", "' + $beteckning + '")'
    $ws.Range("V$rowNum").Formula = '=HYPERLINK("https://klasma.github.io/Logging_JOKKMOKK/klagomål/' + $beteckning + '.docx", "' + $beteckning + '")'
    $ws.Range("W$rowNum").Formula = '=HYPERLINK("https://klasma.github.io/Logging_JOKKMOKK/klagomålsmail/' + $beteckning + '.docx", "' + $beteckning + '")'
    $ws.Range("X$rowNum").Formula = '=HYPERLINK("https://klasma.github.io/Logging_JOKKMOKK/tillsyn/' + $beteckning + '.docx", "' + $beteckning + '")'
    $ws.Range("Y$rowNum").Formula = '=HYPERLINK("https://klasma.github.io/Logging_JOKKMOKK/tillsynsmail/' + $beteckning + '.docx", "' + $beteckning + '")'
}

Set-BeteckningLinks 25 "A 23018-2023"
Set-BeteckningLinks 26 "A 51520-2020"
Set-BeteckningLinks 27 "A 65342-2020"

# The wrapped R-column text triggers row auto-fit; the source file keeps a
# fixed 15pt custom row height throughout, so restore that explicitly.
$ws.Rows.Item(25).RowHeight = 15
$ws.Rows.Item(26).RowHeight = 15
$ws.Rows.Item(27).RowHeight = 15
